$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6136646270751953
$ws.Range("B1").Value = 0.5036099553108215
$ws.Range("C1").Value = 0.4500625729560852
$ws.Range("D1").Value = 0.4969479143619537
$ws.Range("E1").Value = 0.6217054724693298
